$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(5, 16, 4, 4),
    @(2, 15, 4, 5),
    @(4, 5, 3, 15),
    @(3, 5, 4, 15),
    @(7, 19, 4, 1),
    @(5, 12, 4, 8),
    @(5, 7, 4, 13),
    @(5, 15, 4, 5),
    @(2, 2, 3, 18),
    @(6, 8, 4, 12),
    @(4, 8, 2, 12),
    @(3, 6, 5, 14),
    @(6, 7, 9, 13)
)

$startRow = 976
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$ws.Range("A989").Select()
